$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "-"

$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "[-, 'MCT-2A-Programação de Computadores', -, -]"
$ws.Range("F4").Value = "MCT-2A-Circuitos elétricos 2"

$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "[-, 'MCT-2A-Programação de Computadores', -, -]"
$ws.Range("F6").Value = "MCT-2A-Circuitos elétricos 2"

$ws.Range("C7").Value = "-"
$ws.Range("D7").Value = "[-, 'MCT-2A-Programação de Computadores', -, -]"

$ws.Range("B8").Value = "-"
$ws.Range("C8").Value = "-"
$ws.Range("D8").Value = "[-, 'MCT-2A-Programação de Computadores', -, -]"
